$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the URL text from A2 and apply the built-in "Hyperlink" cell style
# (extended report formatting), which adds an underlined, theme-colored font.
$ws.Range("A2").ClearContents()
$ws.Range("A2").Style = "Hyperlink"

# Move / record the active selection as it was left after the edit
$ws.Range("C11").Select()
